$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title: "ubuntu16.04安装spark-2.3.1" -> "Centos7" + new _GoBack bookmark
#    + "安装spark-2.3.1", where only the "Centos7" run gets an explicit
#    en-US / zh-CN language tag.
#
#    The existing _GoBack bookmark (near the end of the document) is moved
#    here, so first drop it from its old location.
# ---------------------------------------------------------------------------

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$d.Content.Find.Execute("ubuntu16.04安装spark-2.3.1", $true, $false, $false, $false, $false, $true, 1, $false, "Centos7", 2) | Out-Null

# Locate the split point (end of the freshly-written "Centos7") dynamically
# rather than hard-coding its length.
$titleFind = $d.Content
$titleFind.Find.Execute("Centos7", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $titleFind.End

# Split "Centos7" into its own paragraph momentarily so the LanguageID /
# LanguageIDFarEast assignment (which this host applies per-paragraph)
# only touches this run, then rebuild "安装spark-2.3.1" in the now-empty
# second paragraph, and finally merge the break back out so both runs sit
# in one paragraph again (each keeping its own run formatting).
$breakPoint = $d.Range($splitPos, $splitPos)
$breakPoint.InsertParagraphAfter()

$titlePara = $d.Paragraphs(1).Range
$titlePara.LanguageID = "en-US"
$titlePara.LanguageIDFarEast = "zh-CN"

$secondPara = $d.Paragraphs(2).Range
$restInsert = $d.Range($secondPara.Start, $secondPara.Start)
$restInsert.InsertAfter("安装spark-2.3.1")

$markRange = $d.Range($splitPos, $splitPos + 1)
$markRange.Delete()

$bookRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bookRange)

# ---------------------------------------------------------------------------
# 2) "Cd  " + "/usr/spark/hadoop-2.7.6/etc/hadoop/" -> single run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Cd  /usr/spark/hadoop-2.7.6/etc/hadoop/", $true, $false, $false, $false, $false, $true, 1, $false, "Cd  /usr/spark/hadoop-2.7.6/etc/hadoop/", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "然后vim " + "core-site.xml  " -> single run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("然后vim core-site.xml  ", $true, $false, $false, $false, $false, $true, 1, $false, "然后vim core-site.xml  ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "比如，我的是`" + "export JAVA_HOME=/root/java/jdk1.8.0_171" -> single run
# ---------------------------------------------------------------------------
$exportText = "比如，我的是``export JAVA_HOME=/root/java/jdk1.8.0_171"
$d.Content.Find.Execute($exportText, $true, $false, $false, $false, $false, $true, 1, $false, $exportText, 2) | Out-Null

Write-Output "done"
